$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend row 3 with new sentiment-delta / direction columns
$ws.Range("X3").Value = -1.3299870000000169
$ws.Range("Y3").Value = "Down"

# Add new row 4 (latest trading day) mirroring the existing row layout
$ws.Range("A4").Value = 42633.888368055559

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = "Neutral"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = "Random"
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 1.76

$ws.Range("S4").Value = 0.1055
$ws.Range("S4").NumberFormat = $ws.Range("S3").NumberFormat

$ws.Range("T4").Value = -6.67
$ws.Range("U4").Value = 5.83
$ws.Range("V4").Value = "N/A"
$ws.Range("W4").Value = 0
